$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Q1 table (columns F:I, rows 2-13) ---
# Row 7: outlier 6747 replaced by the value that was in G8 (6661)
$ws.Range("G7").Value = 6661
# Row 8: G8 and I8 outliers removed (values folded away / cleared)
$ws.Range("G8").ClearContents()
$ws.Range("I8").ClearContents()

# --- Q2 table (columns A:D, rows 16-27) ---
# Row 20: outlier 6639 replaced by the value that was in D21 (6784)
$ws.Range("D20").Value = 6784
# Row 21: D21 now takes the value that was in D22 (6796)
$ws.Range("D21").Value = 6796
# Row 22: C22 and D22 outliers removed
$ws.Range("C22").ClearContents()
$ws.Range("D22").ClearContents()

# --- Highlight the corrected averages with Excel's built-in "Good" style ---
$ws.Range("C13").Style = "Good"
$ws.Range("G13").Style = "Good"
$ws.Range("B27").Style = "Good"

# --- Update the saved selection ---
[void]$ws.Range("H15").Select()
